$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.234.76'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").Value = '2.282.47'
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.19'
$ws.Range("E5").Value = '  +1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.54'
$ws.Range("E6").Value = '  -2.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.607'
$ws.Range("E9").Value = '  -0.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.82'
$ws.Range("E10").Value = '  +0.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0903'
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.31'
$ws.Range("E12").Value = '  -1.55%  '
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.966'
$ws.Range("E14").Value = '  -1.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.11'
$ws.Range("E15").Value = '  -2.39%  '
$ws.Range("D16").Value = '2.629.99'
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").Value = '2.284.75'
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").Value = '42.224.27'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.36'
$ws.Range("E19").Value = '  -4.53%  '
$ws.Range("E20").Value = '  -0.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.89'
$ws.Range("E21").Value = '  +28.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.63'
$ws.Range("E22").Value = '  +2.88%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.02'
$ws.Range("E23").Value = '  -1.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.65'
$ws.Range("E24").Value = '  -6.03%  '
$ws.Range("E25").Value = '  -3.30%  '
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("E27").Value = '  -0.87%  '
$ws.Range("E28").Value = '  +2.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.43'
$ws.Range("E29").Value = '  -4.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.03'
$ws.Range("E30").Value = '  +7.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '164.14'
$ws.Range("E31").Value = '  -0.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.10'
$ws.Range("E32").Value = '  +3.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0878'
$ws.Range("E33").Value = '  -0.60%  '
$ws.Range("E34").Value = '  +0.85%  '
$ws.Range("E35").Value = '  -13.02%  '
$ws.Range("E36").Value = '  -3.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.60'
$ws.Range("E37").Value = '  -1.33%  '
$ws.Range("E39").Value = '  +2.10%  '
$ws.Range("E40").Value = '  -6.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.51'
$ws.Range("E41").Value = '  +1.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '68.95'
$ws.Range("E42").Value = '  -2.42%  '
$ws.Range("E43").Value = '  -0.77%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.62'
$ws.Range("E45").Value = '  -11.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.27'
$ws.Range("E46").Value = '  +1.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '112.96'
$ws.Range("E47").Value = '  -3.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '80.02'
$ws.Range("E48").Value = '  +2.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.91'
$ws.Range("E49").Value = '  -2.24%  '
$ws.Range("E50").Value = '  -2.19%  '
$ws.Range("D51").Value = '1.590.69'
$ws.Range("E51").Value = '  +1.74%  '
